$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new CvLAC entry (Carlos Julio Cortes) ---
# Write order matters for shared-string table ordering (matches authored file):
# D4 (url) first, then B4 (given names), then C4 (surname), then A4 (document id).
$googleUrl = "https://www.google.com/url?q=http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh%3D0000142522&sa=D&source=hangouts&ust=1521740785710000&usg=AFQjCNGB77q4x9G4ftIe4x9qdZjq8byGLA"
$plainUrl = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0000142522"

$ws.Range("D4").Value = $plainUrl
$ws.Range("B4").Value = "Carlos Julio"
$ws.Range("C4").Value = "Cortes"
$ws.Range("A4").Value = 123

# Hyperlink for the new CvLAC cell (mirrors a Google-redirect link pasted from Docs/Sheets).
$ws.Hyperlinks.Add($ws.Range("D4"), $googleUrl, "", "", $googleUrl)

# Hyperlinks.Add applies the plain underlined-black font; re-apply the real Hyperlink style.
$ws.Range("D4").Style = "Hyperlink"

# --- B3 ("Manuel") formatting: underlined, black text (previously a bare hyperlink look) ---
$ws.Range("B3").Font.Underline = $true
$ws.Range("B3").Font.Color = 0

# --- Selection moved to A5 ---
$ws.Range("A5").Select() | Out-Null
